$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 1658
$ws1.Range("F3").Value = 9510
$ws1.Range("F5").Value = 755
$ws1.Range("F6").Value = 624
$ws1.Range("F7").Value = 226
$ws1.Range("G7").Value = 6.8
$ws1.Range("F10").Value = 74
$ws1.Range("F11").Value = 1420
$ws1.Range("F14").Value = 1484
$ws1.Range("F16").Value = 305
$ws1.Range("F18").Value = 145
$ws1.Range("F19").Value = 86
$ws1.Range("F20").Value = 393
$ws1.Range("F22").Value = 102
$ws1.Range("F23").Value = 24
$ws1.Range("F25").Value = 49
$ws1.Range("F29").Value = 76
$ws1.Range("F30").Value = 606
$ws1.Range("F32").Value = 7
$ws1.Range("F34").Value = 170
$ws1.Range("F36").Value = 184
$ws1.Range("F37").Value = 324
$ws1.Range("F39").Value = 265
$ws1.Range("F40").Value = 622
$ws1.Range("F42").Value = 743
$ws1.Range("F45").Value = 321

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F8").Value = 9
$ws2.Range("F11").Value = 696
$ws2.Range("F15").Value = 88
$ws2.Range("F24").Value = 25
$ws2.Range("F27").Value = 6
$ws2.Range("C29").Value = "上海·【早鸟5折】致敬大师·告别之作《你想活出怎样的人生》·宫崎骏&久石让 经典动漫作品音乐会（取消）"
$ws2.Range("G29").Value = "不可售"
$ws2.Range("F30").Value = 22
$ws2.Range("F31").Value = 124
$ws2.Range("F39").Value = 27

# Sheet 3: 本地生活
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F5").Value = 371
$ws3.Range("F7").Value = 2393
$ws3.Range("F8").Value = 3643
$ws3.Range("F9").Value = 20
$ws3.Range("F11").Value = 95
$ws3.Range("F12").Value = 113

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 1658
$ws4.Range("F4").Value = 9510
$ws4.Range("F5").Value = 371
$ws4.Range("F7").Value = 3643
$ws4.Range("F8").Value = 20
$ws4.Range("F9").Value = 755
$ws4.Range("F10").Value = 95
$ws4.Range("F11").Value = 95
$ws4.Range("F13").Value = 226
$ws4.Range("G13").Value = 6.8
$ws4.Range("F16").Value = 696
$ws4.Range("F17").Value = 1420
$ws4.Range("F19").Value = 113
$ws4.Range("F20").Value = 1484
$ws4.Range("F22").Value = 305
$ws4.Range("F24").Value = 145
$ws4.Range("F26").Value = 102
$ws4.Range("F28").Value = 49
$ws4.Range("F35").Value = 76
$ws4.Range("F36").Value = 606
$ws4.Range("F38").Value = 170
$ws4.Range("F41").Value = 324
$ws4.Range("F44").Value = 623
$ws4.Range("F46").Value = 743
$ws4.Range("F50").Value = 321
